# The edit inserts one new data row into the "Hortaliza ... Ají" table at
# row 106 (pushing the existing rows 106-212 down to 107-213), and fills
# the newly inserted row with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 106; this shifts rows 106:212 down to
# 107:213, growing the sheet from 212 to 213 rows (matches the dimension
# change from A1:R212 to A1:R213 in the target workbook).
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new record.
$ws.Range("A106").Value = 8
$ws.Range("B106").Value = "Terminal La Palmera de La Serena"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 44673
$ws.Range("E106").Value = 4
$ws.Range("F106").Value = 100112021
$ws.Range("G106").Value = "Ají"
$ws.Range("H106").Value = "Inferno"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 400
$ws.Range("K106").Value = 23000
$ws.Range("L106").Value = 24000
$ws.Range("M106").Value = 23500
$ws.Range("N106").Value = "$/caja 15 kilos"
$ws.Range("O106").Value = "Provincia de Limarí"
$ws.Range("P106").Value = 1567
$ws.Range("Q106").Value = 15
$ws.Range("R106").Value = "Hortaliza"
